# Apply cryptos list price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.700.78"
$ws.Range("E2").Value = "  -3.13%  "

$ws.Range("D3").Value = "1.952.28"
$ws.Range("E3").Value = "  -2.57%  "

$ws.Range("D4").Value = "'1.015"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'321.40"

$ws.Range("D6").Value = "'1.013"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4753"

$ws.Range("D8").Value = "'0.4021"
$ws.Range("E8").Value = "  -5.12%  "

$ws.Range("D9").Value = "'54.00"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").Value = "'0.08466"
$ws.Range("E10").Value = "  -5.89%  "

$ws.Range("D11").Value = "'1.058"
$ws.Range("E11").Value = "  -5.50%  "

$ws.Range("D12").Value = "'22.19"
$ws.Range("E12").Value = "  -5.02%  "

$ws.Range("D13").Value = "1.986.86"
$ws.Range("E13").Value = "  -1.32%  "

$ws.Range("D14").Value = "'7.589"
$ws.Range("E14").Value = "  -5.91%  "

$ws.Range("D15").Value = "'6.194"
$ws.Range("E15").Value = "  -4.47%  "

$ws.Range("D16").Value = "'1.016"

$ws.Range("D17").Value = "'0.00001072"
$ws.Range("E17").Value = "  -3.86%  "

$ws.Range("D18").Value = "'88.86"
$ws.Range("E18").Value = "  -5.55%  "

$ws.Range("D19").Value = "'0.06628"

$ws.Range("D20").Value = "'18.59"
$ws.Range("E20").Value = "  -5.93%  "

$ws.Range("D21").Value = "'1.013"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "'5.783"
$ws.Range("E22").Value = "  -2.95%  "

$ws.Range("D23").Value = "28.733.87"
$ws.Range("E23").Value = "  -3.08%  "

$ws.Range("D24").Value = "'11.50"
$ws.Range("E24").Value = "  -4.13%  "

$ws.Range("D25").Value = "'2.290"
$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").Value = "2.214.77"
$ws.Range("E26").Value = "  -1.71%  "

$ws.Range("D27").Value = "'154.62"
$ws.Range("E27").Value = "  -3.14%  "

$ws.Range("D28").Value = "'20.14"
$ws.Range("E28").Value = "  -2.89%  "

$ws.Range("D29").Value = "'5.914"
$ws.Range("E29").Value = "  -8.07%  "

$ws.Range("D30").Value = "'2.144"
$ws.Range("E30").Value = "  -6.68%  "

$ws.Range("D31").Value = "'123.49"
$ws.Range("E31").Value = "  -4.08%  "

$ws.Range("D32").Value = "'1.000"
$ws.Range("E32").Value = "  -5.03%  "

$ws.Range("D33").Value = "'0.09562"
$ws.Range("E33").Value = "  -3.77%  "

$ws.Range("D34").Value = "'5.656"
$ws.Range("E34").Value = "  -3.17%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'3.670"
$ws.Range("E35").Value = "  -3.52%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.430"
$ws.Range("E36").Value = "  -8.72%  "

$ws.Range("E37").Value = "  -5.06%  "

$ws.Range("D38").Value = "'1.265"
$ws.Range("E38").Value = "  -3.30%  "

$ws.Range("D39").Value = "'0.06206"
$ws.Range("E39").Value = "  -2.26%  "

$ws.Range("D40").Value = "'8.716"
$ws.Range("E40").Value = "  -7.15%  "

$ws.Range("D41").Value = "'0.6193"
$ws.Range("E41").Value = "  -5.87%  "

$ws.Range("D42").Value = "'11.04"
$ws.Range("E42").Value = "  -5.60%  "

$ws.Range("D43").Value = "'1.012"
$ws.Range("E43").Value = "  -0.03%  "

$ws.Range("D44").Value = "'0.1914"
$ws.Range("E44").Value = "  -6.80%  "

$ws.Range("E45").Value = "  +2.50%  "

$ws.Range("D46").Value = "'0.5914"
$ws.Range("E46").Value = "  -6.77%  "

$ws.Range("D47").Value = "'12.91"
$ws.Range("E47").Value = "  -4.77%  "

$ws.Range("D48").Value = "'2.066"
$ws.Range("E48").Value = "  -6.25%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.00000000336"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("B50").Value = "PancakeSwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D50").Value = "'3.407"
$ws.Range("E50").Value = "  -3.41%  "

$ws.Range("D51").Value = "'0.06829"
$ws.Range("E51").Value = "  -2.41%  "

